# Slide 4: rearrange the "Lanos" / sales-figure labels for the first table row.
#
# The two labels ("Text 27" = "Lanos", "Text 28" = the units-sold figure) are
# repositioned down to sit under the "Verna" row (which itself moves up),
# the figure's text/size changes from "2,652" to "1700", and both labels are
# moved to the very end of the shape tree (on top of the z-order, after the
# trailing group shape).
#
# Point values below are chosen (float32-precise) so that PowerPoint's
# point->EMU conversion reproduces the exact target EMU offsets/extents.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$textLanos = $s.Shapes.Item(29)   # id 30, "Text 27" -> "Lanos"
$textFigure = $s.Shapes.Item(30)  # id 31, "Text 28" -> sales figure
$textVerna = $s.Shapes.Item(32)   # id 33, "Text 30" -> "Verna"
$textVernaFigure = $s.Shapes.Item(33) # id 34, "Text 31" -> "2,177"

# Move "Verna" label up into the row formerly occupied by "Lanos".
$textVerna.Left = 45
$textVerna.Top = 278.97796630859375

# Move Verna's sales figure up too, with a small horizontal nudge.
$textVernaFigure.Left = 195.10536193847656
$textVernaFigure.Top = 279.90740966796875

# Reposition + resize the "Lanos" label to the new row.
$textLanos.Left = 44.650001525878906
$textLanos.Top = 296.7584533691406
$textLanos.Width = 128.9553680419922
$textLanos.Height = 15.969449043273926

# Update the sales-figure text and reposition/resize it.
$textFigure.TextFrame.TextRange.Text = "1700"
$textFigure.Left = 195.10536193847656
$textFigure.Top = 298.7586669921875
$textFigure.Width = 174.79464721679688
$textFigure.Height = 15.969449043273926

# Send both repositioned labels to the end of the shape tree (top of z-order),
# figure label first, then the "Lanos" text label, matching the target order.
$textFigure.ZOrder(0)
$textLanos.ZOrder(0)
